# "All pre-database steps done" — remove the second test data row (row 5,
# "Test item 2" / TEST_CN_02 / FN_02 / etc.) from the IMM import template.
# The row itself stays in place (dimension/sheetData keep row 5) but all of
# its content is cleared; only the two date-formatted cells (I5, BV5) keep
# their number-format styling as empty cells, same as Excel's
# Edit > Clear > Contents on a populated row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out every value in data row 5 (keeps cell formatting, e.g. the date
# format on I5/BV5, but removes the shared-string/numeric content).
$ws.Range("A5:EF5").ClearContents()

# With the second test item gone, move the active selection back to the
# top-left of the frozen data area (A4) instead of the old W5.
$ws.Range("A4").Select() | Out-Null
